# Fruta / hortaliza, semanal
# Insert a new weekly record as row 7, pushing the existing data rows
# (old rows 7-47) down to rows 8-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").EntireRow.Insert()

$ws.Cells.Item(7, 1).Value = 9
$ws.Cells.Item(7, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(7, 3).Value = 'Metropolitana'
$ws.Cells.Item(7, 4).Value = 45222
$ws.Cells.Item(7, 5).Value = 13
$ws.Cells.Item(7, 6).Value = 'Fruta'
$ws.Cells.Item(7, 7).Value = 100108
$ws.Cells.Item(7, 8).Value = 'Tropicales y subtropicales'
$ws.Cells.Item(7, 9).Value = 100108003
$ws.Cells.Item(7, 10).Value = 'Maracuyá'
$ws.Cells.Item(7, 11).Value = 'Sin especificar'
$ws.Cells.Item(7, 12).Value = 'Primera'
$ws.Cells.Item(7, 13).Value = 52
$ws.Cells.Item(7, 14).Value = 50000
$ws.Cells.Item(7, 15).Value = 50000
$ws.Cells.Item(7, 16).Value = 50000
$ws.Cells.Item(7, 17).Value = '$/caja 18 kilos'
$ws.Cells.Item(7, 18).Value = 'Perú'
$ws.Cells.Item(7, 19).Value = 2778
$ws.Cells.Item(7, 20).Value = 18
